# Update the build timestamp embedded in the workbook's version strings.
# Old timestamp: February 03 2026 17.29.55 EST
# New timestamp: February 03 2026 18.05.36 EST

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
            $cell.Value2 = $val.Replace($oldStamp, $newStamp)
        }
    }
}
